# Insert two new daily-log rows right before row 699, shifting the
# existing rows 699-740 down to 701-742, then populate the two new
# rows with their data. Mirrors what Excel's "Insert Copied/Blank Rows"
# does when a couple of rows get spliced into the middle of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 699 (each Insert() pushes
# everything at/below that row index down by one).
$ws.Rows.Item(699).Insert()
$ws.Rows.Item(699).Insert()

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text storage so date-looking strings like "2026/01/21" are
    # not auto-parsed into a date serial number, matching the sheet's
    # existing plain-text date column. Reset to the Normal cell style
    # afterwards so no stray number-format style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# New row 699: 2026/01/21, 水, 22, 201
Set-TextCell 699 1 "2026/01/21"
$ws.Cells.Item(699, 2).Value = "水"
$ws.Cells.Item(699, 3).Value = 22
$ws.Cells.Item(699, 4).Value = 201

# New row 700: 2026/01/22, 木, 2, 200
Set-TextCell 700 1 "2026/01/22"
$ws.Cells.Item(700, 2).Value = "木"
$ws.Cells.Item(700, 3).Value = 2
$ws.Cells.Item(700, 4).Value = 200
